# Update the stimuli file names to the new "exp0" audio paths.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "audio/exp0/kick_mp1.wav"
$ws.Range("A3").Value = "audio/exp0/kick_mp2.wav"
$ws.Range("A4").Value = "audio/exp0/kick_iso.wav"
$ws.Range("A5").Value = "audio/exp0/snare_mp1.wav"
$ws.Range("A6").Value = "audio/exp0/snare_mp2.wav"
$ws.Range("A7").Value = "audio/exp0/snare_iso.wav"

# Move the active selection to D16, matching the saved view state.
$ws.Range("D16").Select()
